$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Farida%Elshafeey%NULL%1,                  Rana%Magdi%NULL%1,                  Nader%Hindi%NULL%1,                  Mohamed%Elshebiny%NULL%1,                  Nourhan%Farrag%NULL%1,                  Shahd%Mahdy%NULL%1,                  Mohamed%Sabbour%NULL%1,                  Sara%Gebril%NULL%1,                  Mohamed%Nasser%NULL%1,                  Menna%Kamel%NULL%1,                  Abdelrahman%Amir%NULL%1,                  Moataz%Maher Emara%NULL%1,                  Ashraf%Nabhan%anabhan@med.asu.edu.eg%1]"
$ws.Range("E3").Value = "[Na%Li%NULL%0,                  Lefei%Han%NULL%0,                  Min%Peng%2658706528@qq.com%0,                  Yuxia%Lv%NULL%0,                  Yin%Ouyang%NULL%0,                  Kui%Liu%NULL%0,                  Linli%Yue%NULL%0,                  Qiannan%Li%NULL%0,                  Guoqiang%Sun%NULL%0,                  Lin%Chen%NULL%0,                  Lin%Yang%l.yang@polyu.edu.hk%0]"
$ws.Range("E4").Value = "[J.%Juan%NULL%1,                  M. M.%Gil%NULL%1,                  Z.%Rong%NULL%2,                  Z.%Rong%NULL%0,                  Y.%Zhang%NULL%1,                  H.%Yang%liona.poon@cuhk.edu.hk%1,                  L. C.%Poon%liona.poon@cuhk.edu.hk%1]"
$ws.Range("E5").Value = "[Chunchen%Wu%NULL%1,                  Wenzhong%Yang%NULL%1,                  Xiaoxue%Wu%NULL%1,                  Tianzhu%Zhang%NULL%1,                  Yaoyao%Zhao%NULL%1,                  Wei%Ren%ren33107@163.com%1,                  Jianbo%Xia%xjb915@126.com%2,                  Jianbo%Xia%xjb915@126.com%0]"
